$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- values that used to belong to row 3 (swap), with Q/R rounded
$ws.Range("A2").Value = 111962461
$ws.Range("Q2").Value = 799376
$ws.Range("R2").Value = 7230155
$ws.Range("Z2").Value = "14:23"
$ws.Range("AB2").Value = "14:23"

# Row 3 <- values that used to belong to row 2 (swap), with Q/R rounded
$ws.Range("A3").Value = 111962533
$ws.Range("Q3").Value = 799379
$ws.Range("R3").Value = 7230183
$ws.Range("Z3").Value = "14:28"
$ws.Range("AB3").Value = "14:28"
